# Update Leave Card 5/22/2023 1:34 PM
# Fill in additional leave-credit rows on the "2018 LEAVE CREDITS" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Row 19 - Feb 2023 period, 1.25 VL/SL earned
$ws.Range("A19").Value = 44958
$ws.Range("C19").Value = 1.25

# Row 20 - Mar 2023 period, 1.25 VL/SL earned
$ws.Range("A20").Value = 44986
$ws.Range("C20").Value = 1.25

# Row 21 - Apr 2023 period, 1.25 VL/SL earned
$ws.Range("A21").Value = 45017
$ws.Range("C21").Value = 1.25

# Row 22 - May 2023 period, VL leave taken (2 days) with remarks
$ws.Range("A22").Value = 45047
$ws.Range("B22").Value = "VL(2-0-0)"
$ws.Range("D22").Value = 2
$ws.Range("K22").Value = "5/15,16/2023"

# Row 23 - SL leave taken (1 day) with a dated remark
$ws.Range("B23").Value = "SL(1-0-0)"
$ws.Range("H23").Value = 1
$ws.Range("K23").Value = 45058
$ws.Range("K23").NumberFormat = "mm-dd-yy"
